# Replace the numeric month values (1-12) in the "Mes" column of Tabla2
# with Spanish three-letter month abbreviations ("Ene.", "Feb.", ... "Dic.").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 5; $row -le 84; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $monthNum = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNum]
}
